$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.797.68"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.108.11"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.40"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.96"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "3.105.00"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.43"
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  -2.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.24"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").Value = "3.619.93"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "66.730.61"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").Value = "3.105.94"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.36"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "476.86"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.92"
$ws.Range("E23").Value = "  +5.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.45"
$ws.Range("E24").Value = "  +4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.86"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.95"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.42"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.90"
$ws.Range("E30").Value = "  -3.43%  "
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.59"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").Value = "0.0₃0934"
$ws.Range("E34").Value = "  -8.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.84"
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.976"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.31"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.07"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.309"
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.61"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").Value = "2.801.48"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "376.60"
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.54"
$ws.Range("E47").Value = "  -12.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.20"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.80"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("E51").Value = "  -1.84%  "
